# [R2] "Absence of taxis..." paragraph (Functional Requirements, section [G6]):
# add "taxi average waiting time" as one more event that must be notified,
# alongside the absence of available taxis and the reservations overlaps.
$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "overlaps and taxis assigned",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "overlaps, taxi average waiting time and taxi assigned",
    2
) | Out-Null

# Add a new requirement [R4] under goal [G7], right after [R3], stating that
# customers must receive the taxi code to be able to recognize its driver.
$r3 = $d.Paragraphs.Item(37)
$r3.Range.InsertParagraphAfter()
$r4 = $d.Paragraphs.Item(38)
$r4.Range.Text = "[R4] Customers must receive the taxi code in order to be able to recognize its driver."
